$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear G99 (will become a blank Summoner Name cell for the new row 99 data)
$ws.Range("G99").ClearContents()

# Row 96
$ws.Range("A96").Value = 94
$ws.Range("B96").Value = 2.383936531395372
$ws.Range("C96").Value = 3065.4
$ws.Range("D96").Value = 0.01186318291444714
$ws.Range("E96").Value = 14.2
$ws.Range("F96").Value = 318.2
$ws.Range("G96").Value = "19 fotsiny adc"
$ws.Range("H96").Value = "NONE"
$ws.Range("I96").Value = 0.2513900194630072
$ws.Range("J96").Value = 11.8
$ws.Range("K96").Value = 0.0083999895266329

# Row 97
$ws.Range("A97").Value = 97
$ws.Range("B97").Value = 11.47138019017781
$ws.Range("C97").Value = 16820.8
$ws.Range("D97").Value = 0.1148893923255633
$ws.Range("E97").Value = 165
$ws.Range("F97").Value = 373.4
$ws.Range("G97").Value = "LS DUFFY"
$ws.Range("H97").Value = "SOLO"
$ws.Range("I97").Value = 0.2602785146347409
$ws.Range("J97").Value = 15.4
$ws.Range("K97").Value = 0.0102324902907245

# Row 98
$ws.Range("A98").Value = 98
$ws.Range("B98").Value = 6.506973090568204
$ws.Range("C98").Value = 9534.4
$ws.Range("D98").Value = 0.04526183093699399
$ws.Range("E98").Value = 68.8
$ws.Range("F98").Value = 370.2
$ws.Range("G98").Value = "BigFather Rengar"
$ws.Range("H98").Value = "SOLO"
$ws.Range("I98").Value = 0.2520432498535652
$ws.Range("J98").Value = 12.8
$ws.Range("K98").Value = 0.008581153657621576

# Row 99
$ws.Range("A99").Value = 100
$ws.Range("B99").Value = 0
$ws.Range("C99").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("F99").Value = 0
$ws.Range("H99").Value = "SOLO"
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0

# Row 100
$ws.Range("A100").Value = 101
$ws.Range("B100").Value = 1.775319622012229
$ws.Range("C100").Value = 3193.8
$ws.Range("D100").Value = 0.01634241245136187
$ws.Range("E100").Value = 29.4
$ws.Range("F100").Value = 41
$ws.Range("G100").Value = "Portgas D Åce "
$ws.Range("H100").Value = "SOLO"
$ws.Range("I100").Value = 0.02279043913285158
$ws.Range("J100").Value = 4.4
$ws.Range("K100").Value = 0.002445803224013341

# Row 101
$ws.Range("A101").Value = 102
$ws.Range("B101").Value = 11.47138019017781
$ws.Range("C101").Value = 16820.8
$ws.Range("D101").Value = 0.1148893923255633
$ws.Range("E101").Value = 165
$ws.Range("F101").Value = 373.4
$ws.Range("G101").Value = "LS DUFFY"
$ws.Range("H101").Value = "SOLO"
$ws.Range("I101").Value = 0.2602785146347409
$ws.Range("J101").Value = 15.4
$ws.Range("K101").Value = 0.0102324902907245

# Row 102
$ws.Range("A102").Value = 103
$ws.Range("B102").Value = 6.506973090568204
$ws.Range("C102").Value = 9534.4
$ws.Range("D102").Value = 0.04526183093699399
$ws.Range("E102").Value = 68.8
$ws.Range("F102").Value = 370.2
$ws.Range("G102").Value = "BigFather Rengar"
$ws.Range("H102").Value = "SOLO"
$ws.Range("I102").Value = 0.2520432498535652
$ws.Range("J102").Value = 12.8
$ws.Range("K102").Value = 0.008581153657621576

# Row 103
$ws.Range("A103").Value = 104
$ws.Range("B103").Value = 2.383936531395372
$ws.Range("C103").Value = 3065.4
$ws.Range("D103").Value = 0.01186318291444714
$ws.Range("E103").Value = 14.2
$ws.Range("F103").Value = 318.2
$ws.Range("G103").Value = "19 fotsiny adc"
$ws.Range("H103").Value = "NONE"
$ws.Range("I103").Value = 0.2513900194630072
$ws.Range("J103").Value = 11.8
$ws.Range("K103").Value = 0.0083999895266329

# New rows 100-103 in column A need the bold/bordered/centered style used by the rest of column A.
# Copy number formatting from A99 (which already has that style) onto the newly written cells.
$ws.Range("A99").Copy() | Out-Null
$ws.Range("A100:A103").PasteSpecial(-4122)
$excel.CutCopyMode = $false
